# This script applies the cryptos-list price/volume refresh described by the
# commit "Updated cryptos list on Sun Jul 30 07:51:53 UTC 2023 with GitHub Actions".
# All Coin/Link/Price/Volume cells are stored as literal TEXT (not numbers/formulas)
# in the workbook, so every write below forces text storage (NumberFormat "@")
# and then restores the "Normal" style so no stray number-format style sticks to
# the cell (matching the original, unstyled data cells).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellAddr, $text) {
    $rng = $ws.Range($cellAddr)
    $rng.NumberFormat = "@"
    $rng.Value = $text
    $rng.Style = "Normal"
}

# Row 2
Set-TextValue "D2" "29.322.63"
Set-TextValue "E2" "  +0.01%  "

# Row 3
Set-TextValue "D3" "1.877.44"
Set-TextValue "E3" "  +0.27%  "

# Row 4
Set-TextValue "E4" "  +0.08%  "

# Row 5
Set-TextValue "D5" "0.7111"
Set-TextValue "E5" "  -0.13%  "

# Row 6
Set-TextValue "D6" "242.38"
Set-TextValue "E6" "  +0.33%  "

# Row 7
Set-TextValue "D7" "1.001"
Set-TextValue "E7" "  +0.14%  "

# Row 8
Set-TextValue "D8" "0.08002"
Set-TextValue "E8" "  +2.77%  "

# Row 9
Set-TextValue "D9" "0.3153"
Set-TextValue "E9" "  +1.42%  "

# Row 10
Set-TextValue "D10" "24.97"
Set-TextValue "E10" "  -0.43%  "

# Row 11
Set-TextValue "D11" "0.08278"
Set-TextValue "E11" "  -1.47%  "

# Row 12
Set-TextValue "D12" "1.897.01"
Set-TextValue "E12" "  +1.79%  "

# Row 13
Set-TextValue "E13" "  +0.19%  "

# Row 14
Set-TextValue "D14" "94.47"
Set-TextValue "E14" "  +3.74%  "

# Row 15
Set-TextValue "D15" "0.7118"
Set-TextValue "E15" "  +0.05%  "

# Row 16
Set-TextValue "D16" "6.351"
Set-TextValue "E16" "  +4.39%  "

# Row 17
Set-TextValue "D17" "0.000008528"
Set-TextValue "E17" "  +3.73%  "

# Row 18
Set-TextValue "D18" "29.345.64"
Set-TextValue "E18" "  +0.06%  "

# Row 19
Set-TextValue "D19" "244.79"
Set-TextValue "E19" "  +1.88%  "

# Row 20
Set-TextValue "D20" "2.142.28"
Set-TextValue "E20" "  +0.94%  "

# Row 21
Set-TextValue "E21" "  +0.45%  "

# Row 23
Set-TextValue "D23" "7.781"

# Row 24
Set-TextValue "D24" "1.002"
Set-TextValue "E24" "  +0.12%  "

# Row 25
Set-TextValue "D25" "0.1554"
Set-TextValue "E25" "  -2.77%  "

# Row 26
Set-TextValue "D26" "9.052"
Set-TextValue "E26" "  +0.30%  "

# Row 27
Set-TextValue "D27" "162.49"
Set-TextValue "E27" "  -0.24%  "

# Row 28
Set-TextValue "E28" "  +0.18%  "

# Row 29
Set-TextValue "E29" "  -0.22%  "

# Row 30
Set-TextValue "D30" "4.418"
Set-TextValue "E30" "  +0.05%  "

# Row 31
Set-TextValue "D31" "4.315"
Set-TextValue "E31" "  +0.20%  "

# Row 32
Set-TextValue "B32" "Hedera"
Set-TextValue "C32" "https://coinranking.com/coin/jad286TjB+hedera-hbar"
Set-TextValue "D32" "0.05372"
Set-TextValue "E32" "  +1.47%  "

# Row 33
Set-TextValue "B33" "Toncoin"
Set-TextValue "C33" "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
Set-TextValue "D33" "1.176"
Set-TextValue "E33" "  -8.73%  "

# Row 34
Set-TextValue "E34" "  -0.09%  "

# Row 35
Set-TextValue "D35" "0.7657"
Set-TextValue "E35" "  +2.81%  "

# Row 36
Set-TextValue "D36" "1.183"
Set-TextValue "E36" "  +0.64%  "

# Row 37
Set-TextValue "D37" "2.691"
Set-TextValue "E37" "  -0.31%  "

# Row 38
Set-TextValue "E38" "  +0.61%  "

# Row 39
Set-TextValue "D39" "1.258.13"
Set-TextValue "E39" "  +2.53%  "

# Row 40
Set-TextValue "D40" "2.752"
Set-TextValue "E40" "  +0.93%  "

# Row 41
Set-TextValue "D41" "6.514"
Set-TextValue "E41" "  -0.73%  "

# Row 42
Set-TextValue "B42" "TrustWalletToken"
Set-TextValue "C42" "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
Set-TextValue "D42" "0.9148"
Set-TextValue "E42" "  +3.12%  "

# Row 43
Set-TextValue "B43" "Quant"
Set-TextValue "C43" "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
Set-TextValue "D43" "112.96"
Set-TextValue "E43" "  +2.20%  "

# Row 44
Set-TextValue "D44" "74.20"
Set-TextValue "E44" "  +2.17%  "

# Row 45
Set-TextValue "E45" "  +8.21%  "

# Row 46
Set-TextValue "E46" "  +0.13%  "

# Row 47
Set-TextValue "D47" "2.046.71"
Set-TextValue "E47" "  +1.40%  "

# Row 48
Set-TextValue "D48" "0.5220"
Set-TextValue "E48" "  +0.43%  "

# Row 49
Set-TextValue "E49" "  -0.21%  "

# Row 50
Set-TextValue "D50" "9.450"
Set-TextValue "E50" "  +0.49%  "

# Row 51
Set-TextValue "D51" "0.4369"
Set-TextValue "E51" "  +1.22%  "
